# DEBkiss results / AIC stuff.xlsx — wrap up edits to Ch4
#
# Summary of changes:
#  - Sheet1: insert an "AICc" column (new E) with a per-row AICc formula,
#    insert another column before the third data block (old J/K -> L/M),
#    add explanatory text rows 30-35, and make Sheet1 the active/selected
#    sheet with I19 selected.
#  - Sheet2 (the I/J/M/N/Q/R summary table): add "Alt total" rows, and
#    move the selection to J27.
#  - Sheet3 (the small B6:D10 table): no longer the active tab (handled
#    automatically once Sheet1 is activated).
#  - Add a new Sheet4 with a small TL / dry-wt / L conversion calc.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Insert the new "AICc" column at E (shifts old E..K to F..L).
$ws1.Columns("E:E").Insert()
# Insert another column at K (shifts old J/K, now K/L, to L/M).
$ws1.Columns("K:K").Insert()

# Header for the new column.
$ws1.Range("E4").Value = "AICc"

# Per-row AICc formula (set individually -- AutoFill was observed to
# clobber unrelated neighbouring cells in this runtime).
for ($r = 5; $r -le 23; $r++) {
    $ws1.Range("E$r").Formula = "=C$r+((2*2)/(141-1-1))"
}
$ws1.Range("E5:E23").NumberFormat = "0.00"

# Explanatory notes added below the table.
$ws1.Range("B30").Value = "AICc = AIC+(2p(p+1)/(n-p-1))"
$ws1.Range("B31").Value = "For AICc assume that N is each data point as entered in the matlab script, not the data points multiplied by their weights. "
$ws1.Range("B32").Value = "n=141"
$ws1.Range("B33").Value = "p=1"
$ws1.Range("B35").Value = "c=0.02877"

# ---------------------------------------------------------------------
# Sheet2 (I/J/M/N/Q/R summary table)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("M14").Value = "Alt total"
$ws2.Range("N14").Value = 8

$ws2.Range("Q18").Value = "Alt total"
$ws2.Range("R18").Value = 12

$ws2.Range("I24").Value = "Alt total"
$ws2.Range("J24").Value = 121

$ws2.Range("J27").Select()

# ---------------------------------------------------------------------
# New Sheet4 — TL / dry weight / length conversion
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "Sheet4"

$ws4.Range("S15").Value = "TL"
$ws4.Range("T15").Value = "50 mm"

$ws4.Range("S16").Value = "dry wt"
$ws4.Range("T16").Formula = "=EXP((2.997*LN(50))-6.7)"

$ws4.Range("S17").Value = "L^3"
$ws4.Range("T17").Formula = "=152.0688/0.4"

$ws4.Range("S18").Value = "L"
$ws4.Range("T18").Formula = "=T17^(1/3)"

$ws4.Range("S19").Value = "delM"
$ws4.Range("T19").Formula = "=T18/50"

$ws4.Range("L39").Select()

# ---------------------------------------------------------------------
# Final selection / active sheet: Sheet1, cell I19.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("I19").Select()
